$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: change a cell's value while re-typing it (text <-> numeric) by
# borrowing the number-format/font "shape" from a template cell that already
# has the desired look, then writing the real value on top.
# ---------------------------------------------------------------------------
function Set-CellAsTextLike {
    param($target, $template, [string]$text)
    $template.Copy() | Out-Null
    $target.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $template.Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = 0
    $target.Value = $text
}

function Set-CellAsNumberLike {
    param($target, $template, $number)
    $template.Copy() | Out-Null
    $target.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $template.Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = 0
    $target.Value = $number
}

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings)
# "Volume 32   Number  35" -> "...36"
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "36"
$ws.Range("A8").Characters(21, 2).Font.Name = "Andale WT"
$ws.Range("A8").Characters(21, 2).Font.Size = 10

# "Report Covering the Week  8/25/2025  Through  8/31/2025"
#   -> "...9/1/2025  Through  9/7/2025"
$ws.Range("C9").Characters(47, 9).Text = "9/7/2025"
$ws.Range("C9").Characters(27, 9).Text = "9/1/2025"
$ws.Range("C9").Characters(27, 8).Font.Name = "Andale WT"
$ws.Range("C9").Characters(27, 8).Font.Size = 10
$ws.Range("C9").Characters(47, 8).Font.Name = "Andale WT"
$ws.Range("C9").Characters(47, 8).Font.Size = 10

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
Set-CellAsTextLike $ws.Range("C16") $ws.Range("D16") "0"
$ws.Range("F16").Value = 1
$ws.Range("H16").Value = -80
$ws.Range("L16").Value = -64.705882352941
$ws.Range("N16").Value = -96.026490066225

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("F17").Value = 2
$ws.Range("H17").Value = 100
$ws.Range("N17").Value = -54.838709677419

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
Set-CellAsTextLike $ws.Range("C19") $ws.Range("D19") "0"
$ws.Range("F19").Value = 7
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 75
$ws.Range("L19").Value = -5.405405405405
$ws.Range("N19").Value = -73.880597014925

# ---------------------------------------------------------------------------
# Row 21 (TOTAL, bold row)
# ---------------------------------------------------------------------------
Set-CellAsTextLike $ws.Range("C21") $ws.Range("D21") "0"
Set-CellAsTextLike $ws.Range("D21") $ws.Range("D21") "0"
Set-CellAsTextLike $ws.Range("E21") $ws.Range("E21") "***.*"
Set-CellAsNumberLike $ws.Range("F21") $ws.Range("F21") 10
Set-CellAsNumberLike $ws.Range("G21") $ws.Range("G21") 11
$ws.Range("H21").Value = -9.090909090909
Set-CellAsNumberLike $ws.Range("I21") $ws.Range("I21") 57
Set-CellAsNumberLike $ws.Range("J21") $ws.Range("J21") 82
$ws.Range("L21").Value = -12.307692307692
$ws.Range("N21").Value = -83.620689655172

# ---------------------------------------------------------------------------
# Row 24 (Transit)
# ---------------------------------------------------------------------------
Set-CellAsTextLike $ws.Range("C24") $ws.Range("D19") "0"
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = -100
$ws.Range("G24").Value = 7
$ws.Range("H24").Value = -28.571428571428
$ws.Range("J24").Value = 28
$ws.Range("K24").Value = -21.428571428571
$ws.Range("L24").Value = -24.137931034482
$ws.Range("M24").Value = -60

# ---------------------------------------------------------------------------
# Row 26 (Petit Larceny)
# ---------------------------------------------------------------------------
Set-CellAsTextLike $ws.Range("D26") $ws.Range("D19") "0"
Set-CellAsTextLike $ws.Range("E26") $ws.Range("E16") "***.*"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("L26").Value = -41.666666666666

# ---------------------------------------------------------------------------
# Row 28 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = 700
$ws.Range("I28").Value = 23
$ws.Range("K28").Value = 187.5
$ws.Range("L28").Value = 43.75

# ---------------------------------------------------------------------------
# Row 31 (Shooting Inc.)
# ---------------------------------------------------------------------------
Set-CellAsNumberLike $ws.Range("D31") $ws.Range("J31") 1
Set-CellAsNumberLike $ws.Range("E31") $ws.Range("K31") -100
Set-CellAsNumberLike $ws.Range("G31") $ws.Range("J31") 1
Set-CellAsNumberLike $ws.Range("H31") $ws.Range("K31") -100
$ws.Range("J31").Value = 3

# ---------------------------------------------------------------------------
# Row 46 (Historical perspective TOTAL row) - values unchanged, style index
# realignment is an internal workbook-table artifact that does not change the
# visible numbers.
# ---------------------------------------------------------------------------
$ws.Range("C46").Value = 368
$ws.Range("E46").Value = 475
$ws.Range("G46").Value = 187
$ws.Range("I46").Value = 127
$ws.Range("J46").Value = 99

$excel.CutCopyMode = 0
